$d = $word.ActiveDocument

# Locate the "2a)" Heading 1 paragraph - it is immediately preceded by a
# long run of empty "Titulo"-styled paragraphs.
$findRange = $d.Content
[void]$findRange.Find.Execute("2a)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Work out this paragraph's 1-based index by counting the paragraphs that
# precede it.
$precedingCount = $d.Range(0, $findRange.Start).Paragraphs.Count
$headingIndex = $precedingCount + 1

# 1) Remove the last 6 empty "Titulo" paragraphs directly preceding the
#    heading.
$firstEmptyIndex = $headingIndex - 6
$lastEmptyIndex = $headingIndex - 1
$removeRange = $d.Range($d.Paragraphs.Item($firstEmptyIndex).Range.Start, $d.Paragraphs.Item($lastEmptyIndex).Range.End)
$removeRange.Delete()

# 2) Move the "_GoBack" bookmark to the start of the "2a)" heading
#    paragraph (its index shifted down by 6 after the deletion above).
#    Adding a new "_GoBack" bookmark automatically replaces any previous
#    one, matching Word's singleton behaviour for that hidden bookmark
#    (it currently sits on the "2b)" heading).
$newHeadingPara = $d.Paragraphs.Item($headingIndex - 6)
$bmRange = $newHeadingPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
